$d = $word.ActiveDocument

# Replace the three-run paragraph ("En la siguiente imagen..." / "Mostrar" / " Tratamiento...")
# with a single run containing the new paragraph text.
$old = "En la siguiente imagen podemos ver un trozo del report Mostrar Tratamiento en el que podemos ver como esto ralentiza nuestro sistema."
$new = "Estos problemas fueron detectados durante la realización del test de rendimiento correspondiente a la Historia de Usuario 20: Añadir Tratamiento. Dado que cuando añades un tratamiento eres redirigido a la vista de Informe, cada usuario al añadir un informe tarda un poco más que el anterior en acceder a la vista de Informe."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# Insert a new paragraph right after it containing "[Captura de Gatling]".
$p = $d.Paragraphs(3)
$p.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(4)
$p2.Range.Text = "[Captura de Gatling]"
